$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new exam result (row) was uploaded. It becomes the newest entry and is
# inserted as the new row 2, pushing the previous rows 2-6 down to rows 3-7.
# The rank/count column (C) is recomputed (countdown from the new total),
# and the "Fecha" column (D) picks up the next value in sequence. The new
# row itself has no duration/score data yet (E blank, F/G/H/I = 0) -- same
# shape as the original row 2 pattern.

# Capture the existing rows 2-6 (bottom to top) before overwriting, then
# write them back shifted down by one row.
$sourceRows = @(2, 3, 4, 5, 6)
$rowsData = @{}
foreach ($r in $sourceRows) {
    $rowsData[$r] = @{
        A = $ws.Range("A$r").Value()
        B = $ws.Range("B$r").Value()
        D = $ws.Range("D$r").Value()
        E = $ws.Range("E$r").Value()
        F = $ws.Range("F$r").Value()
        G = $ws.Range("G$r").Value()
        H = $ws.Range("H$r").Value()
        I = $ws.Range("I$r").Value()
    }
}

# Write rows 7 down to 3 (old rows 6 down to 2) first, from the bottom up,
# so we never overwrite a source row before it has been read.
for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $data = $rowsData[$r]

    $ws.Range("A$dest").Value = $data.A
    $ws.Range("B$dest").Value = $data.B
    $ws.Range("C$dest").Value = 8 - $dest
    $ws.Range("D$dest").Value = $data.D

    if ($null -eq $data.E) {
        $ws.Range("E$dest").ClearContents()
    } else {
        $ws.Range("E$dest").Value = $data.E
    }

    $ws.Range("F$dest").Value = $data.F
    $ws.Range("G$dest").Value = $data.G
    $ws.Range("H$dest").Value = $data.H
    $ws.Range("I$dest").Value = $data.I
}

# New row 2: the freshly uploaded exam result.
$ws.Range("A2").Value = "Edison"
$ws.Range("B2").Value = "Risaralda"
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = "2025-03-17 13:35:04"
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
